$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto Price (column D) and Volume(1h) (column E) text values.
# D-column prices are leading-apostrophe-quoted so Excel keeps them as
# plain text (matching the source inline-string cells) instead of auto-
# converting dotted/decimal-looking values into numbers.

$ws.Range("D2").Value = "'26.427.64"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "'1.675.02"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").Value = "'221.23"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'0.5360"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").Value = "'0.2679"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.06423"
$ws.Range("D10").Value = "'21.08"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "'0.07852"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "'4.566"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'1.674.61"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "'1.904.45"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'0.5660"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "'0.0₅8217"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'66.58"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'26.498.25"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'4.734"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").Value = "'198.67"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").Value = "'10.37"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").Value = "'6.091"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").Value = "'146.76"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "'0.1238"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'7.285"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "'16.30"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").Value = "'1.516"
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").Value = "'0.05903"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").Value = "'1.290"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "'3.590"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'3.319"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "'0.9748"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").Value = "'2.852"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").Value = "'2.449"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").Value = "'0.5850"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("D39").Value = "'0.01616"
$ws.Range("D40").Value = "'1.082.18"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D41").Value = "'5.935"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "'0.8684"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'104.57"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'1.812.69"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "'58.60"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'0.4406"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").Value = "'8.074"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").Value = "'0.05169"
$ws.Range("E51").Value = "  +0.37%  "
